$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the previous "Week 7" block (rows 23-25) down onto
# the new rows so the new entries look consistent with the rest of the log.
$ws.Range("A23:D23").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B23:D23").Copy()
$ws.Range("B28:D28").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B23:D23").Copy()
$ws.Range("B29:D29").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B23:D23").Copy()
$ws.Range("B30:D30").PasteSpecial(-4122) # xlPasteFormats

$ws.Rows("27").RowHeight = 15.75
$ws.Rows("28").RowHeight = 15.75
$ws.Rows("29").RowHeight = 15.75
$ws.Rows("30").RowHeight = 15.75

# Week 8 entries
$ws.Range("A27").Value = "Week 8"
$ws.Range("B27").Value = "Nov 30th"
$ws.Range("C27").Value = "group meeting for presentation slides"
$ws.Range("D27").Value = 2

$ws.Range("B28").Value = "Dec 3rd"
$ws.Range("C28").Value = "group meeting to delegate final paper tasks, first pass at editing"
$ws.Range("D28").Value = 1

$ws.Range("B29").Value = "Dec 4th"
$ws.Range("C29").Value = "fixed citation numbers, more general editing"
$ws.Range("D29").Value = 2

$ws.Range("B30").Value = "Dec 5th"
$ws.Range("C30").Value = "final editing pass"
$ws.Range("D30").Value = 5

$ws.Range("C31").Select() | Out-Null
